# [NTM:SPACE] Ike + Laythe temperature rebalance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- TEMP TERRA (upper table, rows 1-10): Ike row (row 7) was all "?C" placeholders ---
$ws.Range("H7").Value = "-58℃"
$ws.Range("I7").Value = "-4℃"
$ws.Range("J7").Value = "-58℃"
$ws.Range("K7").Value = "-112℃"

# --- TEMP TERRA (lower table, rows 12-21) ---
# Moho row (row 13): DAWN/DAY/DUSK/NIGHT
$ws.Range("H13").Value = "24℃"
$ws.Range("I13").Value = "35℃"
$ws.Range("K13").Value = "15℃"

# Eve row (row 14)
$ws.Range("H14").Value = "23℃"
$ws.Range("I14").Value = "32℃"
$ws.Range("K14").Value = "12℃"

# Duna/I row (row 17) - only NIGHT changes
$ws.Range("K17").Value = "7℃/-23℃"

# Ike row (row 18) - previously empty placeholders, now filled in
$ws.Range("H18").Value = "-6℃"
$ws.Range("I18").Value = "13℃"
$ws.Range("J18").Value = "-7℃"
$ws.Range("K18").Value = "-14℃"

# Update the active selection to match the saved workbook state
$ws.Range("H17").Select()
